$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 12339632
$ws.Cells.Item(40, 9).Value = 3970228.5
$ws.Cells.Item(40, 10).Value = 50001948
$ws.Cells.Item(40, 11).Value = 3970228.5
$ws.Cells.Item(40, 12).Value = 50001948
$ws.Cells.Item(40, 13).Value = -3970053.5
$ws.Cells.Item(40, 14).Value = -50002298
# Row 62
$ws.Cells.Item(62, 8).Value = 4297.273
$ws.Cells.Item(62, 9).Value = 4376.4736
$ws.Cells.Item(62, 11).Value = 4376.4736
$ws.Cells.Item(62, 13).Value = -3752.4736
# Row 65
$ws.Cells.Item(65, 8).Value = 4297.273
$ws.Cells.Item(65, 9).Value = 4376.4736
$ws.Cells.Item(65, 11).Value = 21882.368
$ws.Cells.Item(65, 13).Value = -18762.368
# Row 116
$ws.Cells.Item(116, 8).Value = 4487.5
$ws.Cells.Item(116, 9).Value = 4428.7144
$ws.Cells.Item(116, 10).Value = 4899
$ws.Cells.Item(116, 11).Value = 4428.7144
$ws.Cells.Item(116, 12).Value = 4899
$ws.Cells.Item(116, 13).Value = -986.7143999999998
$ws.Cells.Item(116, 14).Value = -11783
# Row 132
$ws.Cells.Item(132, 8).Value = 1776.4067
$ws.Cells.Item(132, 9).Value = 1789.7931
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 5369.379300000001
$ws.Cells.Item(132, 12).Value = 3000
$ws.Cells.Item(132, 13).Value = -2839.379300000001
$ws.Cells.Item(132, 14).Value = -8060

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 1234.6786
$ws.Cells.Item(2, 9).Value = 1232.9524
$ws.Cells.Item(2, 10).Value = 1239.8572
$ws.Cells.Item(2, 11).Value = 1232.9524
$ws.Cells.Item(2, 12).Value = 1239.8572
$ws.Cells.Item(2, 13).Value = -1119.9524
$ws.Cells.Item(2, 14).Value = -1465.8572
# Row 32
$ws.Cells.Item(32, 8).Value = 5038.5947
$ws.Cells.Item(32, 9).Value = 3612.1528
$ws.Cells.Item(32, 11).Value = 3612.1528
$ws.Cells.Item(32, 13).Value = -3325.1528
# Row 45
$ws.Cells.Item(45, 8).Value = 1533.84
$ws.Cells.Item(45, 9).Value = 1334.1177
$ws.Cells.Item(45, 11).Value = 1334.1177
$ws.Cells.Item(45, 13).Value = -957.1177
# Row 61
$ws.Cells.Item(61, 8).Value = 2211.9333
$ws.Cells.Item(61, 9).Value = 2178.3
$ws.Cells.Item(61, 10).Value = 2279.2
$ws.Cells.Item(61, 11).Value = 2178.3
$ws.Cells.Item(61, 12).Value = 2279.2
$ws.Cells.Item(61, 13).Value = -1966.3
$ws.Cells.Item(61, 14).Value = -2703.2
# Row 63
$ws.Cells.Item(63, 8).Value = 2262.8333
$ws.Cells.Item(63, 9).Value = 2177.875
$ws.Cells.Item(63, 11).Value = 2177.875
$ws.Cells.Item(63, 13).Value = -1491.875
# Row 66
$ws.Cells.Item(66, 8).Value = 2262.8333
$ws.Cells.Item(66, 9).Value = 2177.875
$ws.Cells.Item(66, 11).Value = 10889.375
$ws.Cells.Item(66, 13).Value = -7457.375
# Row 110
$ws.Cells.Item(110, 8).Value = 2442.647
$ws.Cells.Item(110, 9).Value = 1378.1
$ws.Cells.Item(110, 10).Value = 3963.4285
$ws.Cells.Item(110, 11).Value = 1378.1
$ws.Cells.Item(110, 12).Value = 3963.4285
$ws.Cells.Item(110, 13).Value = 666.9000000000001
$ws.Cells.Item(110, 14).Value = -8053.4285
# Row 116
$ws.Cells.Item(116, 8).Value = 1234.6786
$ws.Cells.Item(116, 9).Value = 1232.9524
$ws.Cells.Item(116, 10).Value = 1239.8572
$ws.Cells.Item(116, 11).Value = 1232.9524
$ws.Cells.Item(116, 12).Value = 1239.8572
$ws.Cells.Item(116, 13).Value = 1061.0476
$ws.Cells.Item(116, 14).Value = -5827.8572
# Row 136
$ws.Cells.Item(136, 8).Value = 2211.9333
$ws.Cells.Item(136, 9).Value = 2178.3
$ws.Cells.Item(136, 10).Value = 2279.2
$ws.Cells.Item(136, 11).Value = 6534.900000000001
$ws.Cells.Item(136, 12).Value = 6837.599999999999
$ws.Cells.Item(136, 13).Value = -3984.900000000001
$ws.Cells.Item(136, 14).Value = -11937.6

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 1234.6786
$ws.Cells.Item(3, 9).Value = 1232.9524
$ws.Cells.Item(3, 10).Value = 1239.8572
$ws.Cells.Item(3, 11).Value = 1232.9524
$ws.Cells.Item(3, 12).Value = 1239.8572
$ws.Cells.Item(3, 13).Value = -1118.9524
$ws.Cells.Item(3, 14).Value = -1467.8572
# Row 134
$ws.Cells.Item(134, 8).Value = 1487.4524
$ws.Cells.Item(134, 9).Value = 1495
$ws.Cells.Item(134, 11).Value = 4485
$ws.Cells.Item(134, 13).Value = -1950

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 322723.34
$ws.Cells.Item(31, 9).Value = 4591.225
$ws.Cells.Item(31, 11).Value = 4591.225
$ws.Cells.Item(31, 13).Value = -4296.225
# Row 34
$ws.Cells.Item(34, 8).Value = 322723.34
$ws.Cells.Item(34, 9).Value = 4591.225
$ws.Cells.Item(34, 11).Value = 4591.225
$ws.Cells.Item(34, 13).Value = -4389.225
# Row 94
$ws.Cells.Item(94, 8).Value = 3104
$ws.Cells.Item(94, 10).Value = 4914
$ws.Cells.Item(94, 12).Value = 4914
$ws.Cells.Item(94, 14).Value = -5816
# Row 122
$ws.Cells.Item(122, 8).Value = 1066.0769
$ws.Cells.Item(122, 10).Value = 1176.7142
$ws.Cells.Item(122, 12).Value = 3530.1426
$ws.Cells.Item(122, 14).Value = -8430.142599999999
# Row 134
$ws.Cells.Item(134, 8).Value = 2512.8
$ws.Cells.Item(134, 9).Value = 2414.2222
$ws.Cells.Item(134, 11).Value = 7242.6666
$ws.Cells.Item(134, 13).Value = -4707.6666

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Cells.Item(38, 8).Value = 202
$ws.Cells.Item(38, 9).Value = 300
$ws.Cells.Item(38, 11).Value = 900
$ws.Cells.Item(38, 13).Value = -553
# Row 108
$ws.Cells.Item(108, 8).Value = 41.333332
$ws.Cells.Item(108, 9).Value = 41.333332
$ws.Cells.Item(108, 11).Value = 123.999996
$ws.Cells.Item(108, 13).Value = 2756.000004
# Row 109
$ws.Cells.Item(109, 8).Value = 1659.4
$ws.Cells.Item(109, 9).Value = 1659.4
$ws.Cells.Item(109, 11).Value = 4978.200000000001
$ws.Cells.Item(109, 13).Value = -3938.200000000001
# Row 121
$ws.Cells.Item(121, 8).Value = 83877.664
$ws.Cells.Item(121, 9).Value = 167023.33
$ws.Cells.Item(121, 10).Value = 732
$ws.Cells.Item(121, 11).Value = 501069.99
$ws.Cells.Item(121, 12).Value = 2196
$ws.Cells.Item(121, 13).Value = -499759.99
$ws.Cells.Item(121, 14).Value = -4816
# Row 122
$ws.Cells.Item(122, 8).Value = 914.4286
$ws.Cells.Item(122, 9).Value = 950
$ws.Cells.Item(122, 10).Value = 900.2
$ws.Cells.Item(122, 11).Value = 8550
$ws.Cells.Item(122, 12).Value = 8101.8
$ws.Cells.Item(122, 13).Value = -6100
$ws.Cells.Item(122, 14).Value = -13001.8
# Row 129
$ws.Cells.Item(129, 8).Value = 3367.6
$ws.Cells.Item(129, 10).Value = 3009.9333
$ws.Cells.Item(129, 12).Value = 9029.7999
$ws.Cells.Item(129, 14).Value = -19029.7999
# Row 131
$ws.Cells.Item(131, 8).Value = 1780.0344
$ws.Cells.Item(131, 10).Value = 2044.7
$ws.Cells.Item(131, 12).Value = 6134.1
$ws.Cells.Item(131, 14).Value = -16214.1

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 337.83334
$ws.Cells.Item(2, 9).Value = 367
$ws.Cells.Item(2, 11).Value = 367
$ws.Cells.Item(2, 13).Value = -254
# Row 122
$ws.Cells.Item(122, 8).Value = 2500.5715
$ws.Cells.Item(122, 9).Value = 1325.6
$ws.Cells.Item(122, 10).Value = 3568.7273
$ws.Cells.Item(122, 11).Value = 3976.8
$ws.Cells.Item(122, 12).Value = 10706.1819
$ws.Cells.Item(122, 13).Value = -1526.8
$ws.Cells.Item(122, 14).Value = -15606.1819
# Row 132
$ws.Cells.Item(132, 8).Value = 5634.433
$ws.Cells.Item(132, 9).Value = 4853.826
$ws.Cells.Item(132, 11).Value = 14561.478
$ws.Cells.Item(132, 13).Value = -12031.478

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 40006910
$ws.Cells.Item(7, 9).Value = 90913016
$ws.Cells.Item(7, 10).Value = 9258.5
$ws.Cells.Item(7, 11).Value = 90913016
$ws.Cells.Item(7, 12).Value = 9258.5
$ws.Cells.Item(7, 13).Value = -90912904
$ws.Cells.Item(7, 14).Value = -9482.5
# Row 122
$ws.Cells.Item(122, 8).Value = 5575.8647
$ws.Cells.Item(122, 9).Value = 3573.8948
$ws.Cells.Item(122, 10).Value = 7689.0557
$ws.Cells.Item(122, 11).Value = 10721.6844
$ws.Cells.Item(122, 12).Value = 23067.1671
$ws.Cells.Item(122, 13).Value = -8271.6844
$ws.Cells.Item(122, 14).Value = -27967.1671
# Row 126
$ws.Cells.Item(126, 8).Value = 40006910
$ws.Cells.Item(126, 9).Value = 90913016
$ws.Cells.Item(126, 10).Value = 9258.5
$ws.Cells.Item(126, 11).Value = 272739048
$ws.Cells.Item(126, 12).Value = 27775.5
$ws.Cells.Item(126, 13).Value = -272736578
$ws.Cells.Item(126, 14).Value = -32715.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 49
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 14).Value = $null
# Row 132
$ws.Cells.Item(132, 8).Value = 2631.8
$ws.Cells.Item(132, 9).Value = 2320.5881
$ws.Cells.Item(132, 10).Value = 2820.75
$ws.Cells.Item(132, 11).Value = 6961.7643
$ws.Cells.Item(132, 12).Value = 8462.25
$ws.Cells.Item(132, 13).Value = -4431.7643
$ws.Cells.Item(132, 14).Value = -13522.25
